# Append: 2026-02-10 13:22 JST
# The scraper re-ran and refreshed the "取得日時" (fetched-at) timestamp
# for every already-recorded row on the listing sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2026-02-10 13:07:34"
$newTimestamp = "2026-02-10 13:22:43"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
